$d = $word.ActiveDocument

function Replace-In-Range($range, $old, $new) {
    $f = $range.Find
    $f.ClearFormatting()
    $found = $f.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false)
    if (-not $found) {
        throw "Text not found: $old"
    }
    $range.Text = $new
}

function Replace-Whole-Doc($old, $new) {
    $rng = $d.Content
    Replace-In-Range $rng $old $new
}

# --- Title ---
Replace-Whole-Doc "Cosmic Marvels: Unveiling the Enigma of Dark Matter" "Chemistry: The Science of Change"

# --- Author (merges 3 runs "Dr" + "." + " Riya Sharma" into one run) ---
Replace-Whole-Doc "Dr. Riya Sharma" "Professor John Maxwell"

# --- Email address (two separate runs around a "." run that is left alone) ---
Replace-Whole-Doc "riyasharma@astronews" "jmaxwell@schoolmail"
Replace-In-Range ($d.Paragraphs(3).Range) "com" "edu"

# --- Body paragraph 1: sentence-by-sentence replacements ---
$p5 = $d.Paragraphs(5).Range

Replace-In-Range $p5 "Galaxies, the celestial spectacles that adorn our night sky, have captivated the hearts and minds of astronomers for centuries" "In the realm of natural sciences, Chemistry stands as a captivating and ever-evolving discipline that delves into the composition, structure, properties, and behavior of matter"

Replace-In-Range $p5 "Within these vast cosmic tapestries, a perplexing mystery unfurls--a hidden mass that exerts a profound gravitational influence on the structure and dynamics of our universe" "It encapsulates the study of substances, their interactions with each other, and the energy changes that accompany these interactions"

Replace-In-Range $p5 "This elusive substance, known as dark matter, has baffled scientists for decades, prompting a quest to unravel its enigmatic nature" "Chemistry's vast scope encompasses the investigation of materials at the atomic and molecular levels, including their synthesis, reactivity, and various applications in diverse fields"

Replace-In-Range $p5 "Discoveries in astrophysics, cosmology, and particle physics have shed light on the existence and properties of dark matter" "From the intricate molecular dance that orchestrates biological processes to the marvels of material science and technological advancements, Chemistry plays a pivotal role in shaping our understanding of the world around us"

Replace-In-Range $p5 "Intriguing theories propose exotic particles or modifications to the laws of gravity to account for this mysterious entity" "Its profound influence extends far beyond the laboratory, impacting industries, medicine, agriculture, and countless aspects of our daily lives"

Replace-In-Range $p5 "Despite its elusive nature, the gravitational effects of dark matter are evident throughout the universe" "Through experimentation, observation, and an insatiable quest for knowledge, chemists unlock the secrets of matter"

Replace-In-Range $p5 "Observations of galaxy rotation curves, gravitational lensing, and the behavior of galaxy clusters all point to the presence of a vast reservoir of unseen mass" "They unravel the molecular mechanisms responsible for life's symphony, devise synthetic routes to novel compounds with tailored properties, and explore the intricate interactions between matter and energy"

Replace-In-Range $p5 "Dark matter constitutes approximately 27% of the universe's energy-mass budget, dwarfing the contribution of visible matter" "Chemistry's journey of discovery not only expands our comprehension of the universe but also impinges upon our daily routines"

# This replacement also swallows the following "." run and the whole next sentence run,
# leaving the sentence's trailing "." run (originally terminating the swallowed sentence) intact.
Replace-In-Range $p5 "Unraveling the enigma of dark matter holds the key to comprehending the universe's large-scale structure, the formation and evolution of galaxies, and the nature of gravity itself. The search for dark matter particles at underground laboratories, the analysis of cosmic microwave background radiation, and the exploration of alternative theories of gravity are among the ongoing endeavors to illuminate this cosmic mystery" "The advances in Chemistry have brought forth transformative technologies, such as versatile plastics, potent pharmaceuticals, sustainable energy sources, and cutting-edge materials that revolutionize industries"

Replace-In-Range $p5 "While its true identity remains shrouded in mystery, the study of dark matter has opened up new avenues of exploration in physics, pushing the boundaries of our understanding of the fundamental forces that govern the universe" "Chemistry is more than just an academic pursuit; it is a dynamic and vibrant narrative of scientific exploration, innovation, and societal impact"

Replace-In-Range $p5 "It is a testament to the enduring human curiosity and the relentless pursuit of knowledge that scientists continue to delve into the enigma of dark matter, hoping to unravel its secrets and illuminate the cosmos" "As we delve deeper into the complexities of matter, Chemistry continues to unveil hidden truths, illuminating our comprehension of the microcosm and shaping the course of human progress"

# Append four new runs at the very end of paragraph 5 (after the last, still-unchanged "." run)
$endOfP5 = $d.Paragraphs(5).Range
$endOfP5.SetRange($endOfP5.End - 1, $endOfP5.End - 1)
$endOfP5.InsertAfter(".")
$endOfP5.SetRange($endOfP5.End, $endOfP5.End)
$endOfP5.InsertAfter(" Its significance lies not only in unraveling the fundamental principles that govern the interactions of substances but also in harnessing this knowledge to address global challenges, improve human health, and advance technological frontiers")
$endOfP5.SetRange($endOfP5.End, $endOfP5.End)
$endOfP5.InsertAfter(".")
$endOfP5.SetRange($endOfP5.End, $endOfP5.End)
$endOfP5.InsertAfter(" The pursuit of Chemistry is not merely a journey of academic intrigue; it is an odyssey of exploration, discovery, and transformative change")

Write-Output "body paragraph done"

# --- Summary heading / content ---
$p7 = $d.Paragraphs(7).Range

Replace-In-Range $p7 "Dark matter remains one of the most profound mysteries in modern physics" "In this comprehensive essay, I have endeavored to capture the captivating essence of Chemistry, a science that interrogates the nature of matter, energy, and their intricate interplay"

Replace-In-Range $p7 "Its gravitational influence is evident throughout the universe, yet its true nature remains elusive" "Exploring its diverse facets, from the study of molecular structures to the development of innovative materials, I have elucidated the profound impact Chemistry has on our lives and the world around us"

Replace-In-Range $p7 "While its identity is unknown, the search for dark matter particles and the exploration of alternative theories of gravity are ongoing endeavors" "As we continue to unravel the enigmas of matter, Chemistry stands poised to illuminate new pathways towards scientific breakthroughs and societal progress"

Replace-In-Range $p7 "Unraveling the enigma of dark matter holds the key to understanding the universe's structure, the formation and evolution of galaxies, and the nature of gravity itself. The pursuit of this cosmic mystery continues to drive scientific exploration and push the boundaries of our knowledge" "Its enduring legacy lies in its ability to not only unravel the intricate workings of the universe but also in empowering us to harness its secrets for the betterment of humankind"

Write-Output "summary done"
